# Quarto site content update: add "Exercise time" / "Study time" variables
# to the example dataset (Data sheet) and document them in the Codebook sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Data")
$ws2 = $wb.Worksheets.Item("Codebook")

# --- Data sheet: two new columns, D (Exercise time) and E (Study time) ---
$ws1.Range("D1").Value = "Exercise time"
$ws1.Range("D1").Font.Bold = $true
$ws1.Range("E1").Value = "Study time"
$ws1.Range("E1").Font.Bold = $true

# column widths for the Data sheet (A..E)
$ws1.Columns.Item(1).ColumnWidth = 17.59
$ws1.Columns.Item(2).ColumnWidth = 23.09
$ws1.Columns.Item(3).ColumnWidth = 15.422
$ws1.Columns.Item(4).ColumnWidth = 13.2541
$ws1.Columns.Item(5).ColumnWidth = 16.422

# --- Codebook sheet: document the two new variables ---
$ws2.Range("A5").Value = "Excersise time "
$ws2.Range("B5").Value = "total time spent exercising weekly (minutes)"
$ws2.Range("C5").Value = "numeric value >0 or NA"

$ws2.Range("A6").Value = "Study time"
$ws2.Range("B6").Value = "preffered study time(morning/afternoon/night)"
$ws2.Range("C6").Value = "M/A/N"

# column widths for the Codebook sheet (A..C)
$ws2.Columns.Item(1).ColumnWidth = 13.422
$ws2.Columns.Item(2).ColumnWidth = 38.422
$ws2.Columns.Item(3).ColumnWidth = 20.922

# --- Final on-screen selection / active sheet, matching the saved UI state ---
$ws1.Columns.Item(4).Select() | Out-Null
$ws2.Range("C8").Select() | Out-Null
$ws2.Activate()
